# Update NATMI TPM-derived values on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new value, matching the recomputed NATMI metrics
# after the ligand (Pla2g10 @ MuSCs) / receptor (Pla2r1 @ MuSCs & ECs) TPM update.

# Row 2 (FAPs -> MuSCs)
$ws.Range("I2").Value = 0.4698873783336738
$ws.Range("J2").Value = 0.4698873783336738
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4009296666666667
$ws.Range("N2").Value = 1.202789
$ws.Range("O2").Value = 0.01991582391091379
$ws.Range("P2").Value = 0.01991582391091379
$ws.Range("Q2").Value = 0.028587888952
$ws.Range("R2").Value = 0.257291000568
$ws.Range("S2").Value = 0.009358194284854373
$ws.Range("T2").Value = 0.009358194284854375

# Row 3 (FAPs -> FAPs)
$ws.Range("I3").Value = 0.4698873783336738
$ws.Range("J3").Value = 0.4698873783336738
$ws.Range("O3").Value = 0.4225534031856171
$ws.Range("P3").Value = 0.4225534031856172
$ws.Range("S3").Value = 0.1985525108288615
$ws.Range("T3").Value = 0.1985525108288615

# Row 4 (FAPs -> ECs)
$ws.Range("I4").Value = 0.4698873783336738
$ws.Range("J4").Value = 0.4698873783336738
$ws.Range("M4").Value = 11.22377
$ws.Range("N4").Value = 33.67131000000001
$ws.Range("O4").Value = 0.5575307729034691
$ws.Range("P4").Value = 0.5575307729034691
$ws.Range("Q4").Value = 0.8002996960800001
$ws.Range("R4").Value = 7.202697264720001
$ws.Range("S4").Value = 0.2619766732199579
$ws.Range("T4").Value = 0.2619766732199579

# Row 5 (MuSCs -> ECs)
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.080443
$ws.Range("H5").Value = 0.241329
$ws.Range("I5").Value = 0.5301126216663262
$ws.Range("J5").Value = 0.5301126216663261
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4009296666666667
$ws.Range("N5").Value = 1.202789
$ws.Range("O5").Value = 0.01991582391091379
$ws.Range("P5").Value = 0.01991582391091379
$ws.Range("Q5").Value = 0.03225198517566667
$ws.Range("R5").Value = 0.290267866581
$ws.Range("S5").Value = 0.01055762962605942
$ws.Range("T5").Value = 0.01055762962605941

# Row 6 (MuSCs -> FAPs)
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.080443
$ws.Range("H6").Value = 0.241329
$ws.Range("I6").Value = 0.5301126216663262
$ws.Range("J6").Value = 0.5301126216663261
$ws.Range("O6").Value = 0.4225534031856171
$ws.Range("P6").Value = 0.4225534031856172
$ws.Range("Q6").Value = 0.684289344816
$ws.Range("R6").Value = 6.158604103344
$ws.Range("S6").Value = 0.2240008923567557
$ws.Range("T6").Value = 0.2240008923567556

# Row 7 (MuSCs -> ECs)
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.080443
$ws.Range("H7").Value = 0.241329
$ws.Range("I7").Value = 0.5301126216663262
$ws.Range("J7").Value = 0.5301126216663261
$ws.Range("M7").Value = 11.22377
$ws.Range("N7").Value = 33.67131000000001
$ws.Range("O7").Value = 0.5575307729034691
$ws.Range("P7").Value = 0.5575307729034691
$ws.Range("Q7").Value = 0.9028737301100002
$ws.Range("R7").Value = 8.125863570990001
$ws.Range("S7").Value = 0.2955540996835112
$ws.Range("T7").Value = 0.2955540996835111
